# Apply cryptocurrency price/volume updates from the Feb 9 2024 GitHub Actions refresh run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "47.330.67"
$ws.Range("E2").Value = "  +5.79%  "

$ws.Range("D3").Value = "2.515.28"
$ws.Range("E3").Value = "  +3.94%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D6").Value = "'105.51"
$ws.Range("E6").Value = "  +3.82%  "

$ws.Range("D7").Value = "'0.524"
$ws.Range("E7").Value = "  +2.15%  "

$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").Value = "'0.541"
$ws.Range("E9").Value = "  +2.43%  "

$ws.Range("D10").Value = "'37.11"
$ws.Range("E10").Value = "  +5.04%  "

$ws.Range("E11").Value = "  +2.64%  "

$ws.Range("E12").Value = "  +0.98%  "

$ws.Range("E13").Value = "  -0.43%  "

$ws.Range("E14").Value = "  +4.61%  "

$ws.Range("D15").Value = "2.911.34"
$ws.Range("E15").Value = "  +4.08%  "

$ws.Range("D16").Value = "2.540.98"
$ws.Range("E16").Value = "  +4.00%  "

$ws.Range("D17").Value = "'0.848"
$ws.Range("E17").Value = "  +2.61%  "

$ws.Range("D18").Value = "47.250.01"
$ws.Range("E18").Value = "  +6.01%  "

$ws.Range("D19").Value = "'12.80"
$ws.Range("E19").Value = "  +4.73%  "

$ws.Range("D20").Value = "'6.57"
$ws.Range("E20").Value = "  +3.62%  "

$ws.Range("E21").Value = "  +2.83%  "

$ws.Range("D22").Value = "'70.97"
$ws.Range("E22").Value = "  +3.53%  "

$ws.Range("D23").Value = "'252.19"
$ws.Range("E23").Value = "  +3.91%  "

$ws.Range("D24").Value = "'2.39"
$ws.Range("E24").Value = "  +5.46%  "

$ws.Range("E25").Value = "  +3.39%  "

$ws.Range("D26").Value = "'26.54"
$ws.Range("E26").Value = "  +5.53%  "

$ws.Range("E27").Value = "  -0.03%  "

$ws.Range("E28").Value = "  +5.46%  "

$ws.Range("E29").Value = "  -3.47%  "

$ws.Range("D30").Value = "'35.19"
$ws.Range("E30").Value = "  +5.52%  "

$ws.Range("E31").Value = "  +7.97%  "

$ws.Range("D32").Value = "'49.65"
$ws.Range("E32").Value = "  +2.86%  "

$ws.Range("D33").Value = "'19.86"
$ws.Range("E33").Value = "  +1.76%  "

$ws.Range("E34").Value = "  +3.02%  "

$ws.Range("D35").Value = "'0.0784"
$ws.Range("E35").Value = "  +2.98%  "

$ws.Range("E36").Value = "  +0.03%  "

$ws.Range("E37").Value = "  +4.07%  "

$ws.Range("E38").Value = "  +4.04%  "

$ws.Range("E39").Value = "  +4.37%  "

$ws.Range("D40").Value = "'123.75"
$ws.Range("E40").Value = "  -1.89%  "

$ws.Range("E41").Value = "  +2.34%  "

$ws.Range("E42").Value = "  +2.63%  "

$ws.Range("D43").Value = "'21.35"
$ws.Range("E43").Value = "  +1.93%  "

$ws.Range("E44").Value = "  +3.24%  "

$ws.Range("D45").Value = "1.981.10"
$ws.Range("E45").Value = "  +2.33%  "

$ws.Range("E46").Value = "  +3.63%  "

$ws.Range("E47").Value = "  +1.03%  "

$ws.Range("E48").Value = "  +2.92%  "

$ws.Range("D49").Value = "'9.12"
$ws.Range("E49").Value = "  +0.00%  "

$ws.Range("E50").Value = "  +17.37%  "

$ws.Range("D51").Value = "'79.76"
$ws.Range("E51").Value = "  +5.04%  "
